$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "before" sheet has category-header rows (sexo, cor ou raça, grupos de
# idade, nível de instrução, classes de rendimento...) plus trailing
# footnote rows (sem rendimento a menos, fonte:..., (1) inclusive...) that
# need to be removed entirely. Deleting each row shifts the data below it
# up, which also re-compacts the numeric data into the now-adjacent rows
# and drops the now-unused shared strings automatically on save.
#
# Original row numbers (top to bottom) to remove:
#   5  -> "sexo"
#   8  -> "cor ou raça"
#   13 -> "grupos de idade"
#   19 -> "nível de instrução"
#   27 -> "classes de rendimento mensal domiciliar per capita"
#   29 -> "sem rendimento a menos "
#   35 -> "fonte: ibge, ..."
#   36 -> "(1) inclusive ..."
#
# Because each delete shifts everything below up by one row, the row
# positions to select at each step (applied in order) are:
$ws.Rows("5:5").Delete()
$ws.Rows("7:7").Delete()
$ws.Rows("11:11").Delete()
$ws.Rows("16:16").Delete()
$ws.Rows("23:23").Delete()
$ws.Rows("24:24").Delete()
$ws.Rows("29:29").Delete()
$ws.Rows("29:29").Delete()
